$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-08-07 Thursday" "2025-08-08 Friday"

Replace-Text "714÷4=" "949÷2="
Replace-Text "103÷8=" "386÷5="
Replace-Text "570÷5=" "114÷4="
Replace-Text "877÷9=" "554÷8="
Replace-Text "369÷2=" "886÷9="
Replace-Text "111÷2=" "465÷9="
Replace-Text "602÷9=" "668÷8="
Replace-Text "665÷5=" "260÷9="
Replace-Text "944÷7=" "355÷3="
Replace-Text "673÷6=" "752÷2="
Replace-Text "190÷9=" "779÷3="
Replace-Text "275÷9=" "257÷8="
Replace-Text "613÷7=" "684÷8="
Replace-Text "790÷7=" "165÷7="
Replace-Text "491÷6=" "214÷2="
Replace-Text "943÷7=" "831÷4="
Replace-Text "754÷7=" "936÷9="
Replace-Text "450÷9=" "357÷4="
Replace-Text "510÷3=" "453÷4="
Replace-Text "308÷9=" "155÷2="
Replace-Text "203÷7=" "461÷2="
Replace-Text "715÷4=" "224÷8="
Replace-Text "898÷3=" "842÷3="
Replace-Text "288÷7=" "408÷6="
Replace-Text "435÷5=" "184÷7="
